$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH119-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21D | GRAP COUNT NUMER: NONE"

$ws.Range("A1:H2").Select()
$ws.Range("A5").Select()
